$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.192.72'
$ws.Range("E2").Value = '  +0.18%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.852.23'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.6979'
$ws.Range("E5").Value = '  +1.73%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '236.94'
$ws.Range("E6").Value = '  -0.35%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.001'
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07880'
$ws.Range("E8").Value = '  +1.40%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3015'
$ws.Range("E9").Value = '  -0.70%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.81'
$ws.Range("E10").Value = '  +2.80%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08134'
$ws.Range("E11").Value = '  +0.64%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.888.61'
$ws.Range("E12").Value = '  +2.06%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.179'
$ws.Range("E13").Value = '  -0.35%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7041'
$ws.Range("E14").Value = '  -2.41%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.26'
$ws.Range("E15").Value = '  -0.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.188.78'
$ws.Range("E16").Value = '  +0.20%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.786'
$ws.Range("E17").Value = '  +0.92%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007817'
$ws.Range("E18").Value = '  +0.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.18'
$ws.Range("E19").Value = '  -0.70%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '235.24'
$ws.Range("E20").Value = '  +0.59%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.000'
$ws.Range("E21").Value = '  +0.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.093.78'
$ws.Range("E22").Value = '  -0.62%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.485'
$ws.Range("E24").Value = '  +0.22%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '162.20'
$ws.Range("E25").Value = '  +0.33%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.872'
$ws.Range("E26").Value = '  -1.18%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1411'
$ws.Range("E27").Value = '  -0.95%  '
$ws.Range("E28").Value = '  +0.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.902'
$ws.Range("E29").Value = '  -2.52%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.399'
$ws.Range("E30").Value = '  -0.67%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.472'
$ws.Range("E31").Value = '  -0.49%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.264'
$ws.Range("E32").Value = '  -5.44%  '
$ws.Range("E33").Value = '  -0.04%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05149'
$ws.Range("E34").Value = '  -0.91%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.164'
$ws.Range("E35").Value = '  -0.97%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7054'
$ws.Range("E36").Value = '  +0.24%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9961'
$ws.Range("E37").Value = '  -1.47%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.678'
$ws.Range("E38").Value = '  +0.65%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01843'
$ws.Range("E39").Value = '  -0.34%  '
$ws.Range("E40").Value = '  +0.43%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.147.41'
$ws.Range("E41").Value = '  +3.74%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9182'
$ws.Range("E42").Value = '  -2.40%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.935'
$ws.Range("E43").Value = '  +0.33%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4226'
$ws.Range("E44").Value = '  -1.19%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '69.91'
$ws.Range("E45").Value = '  -0.70%  '
$ws.Range("E46").Value = '  +0.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.80'
$ws.Range("E47").Value = '  +0.38%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5295'
$ws.Range("E48").Value = '  -2.83%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.734'
$ws.Range("E49").Value = '  -3.29%  '
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.993.14'
$ws.Range("E50").Value = '  -0.44%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.132'
$ws.Range("E51").Value = '  -0.16%  '
